# Update weeks/milestones schedule: shift all dates forward by 5 days,
# add milestone text for weeks 1-3 and 12, wrap text on the milestone
# column, increase row heights, and update the saved selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New date values (serial numbers), one per data row 2..13
$dates = @(42408, 42415, 42423, 42429, 42436, 42443, 42450, 42457, 42464, 42471, 42478, 42485)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

# Milestone text for weeks 1 (row2), 2 (row3), 3 (row4) and week 12 (row13)
$wrapUp = "Wrap up for final presentation"
$pmText = "• Project Manager: 1.Jake 2.Colton: Make sure people know what they’re doing.  Come up with weekly sprint objectives.  Have a rough timeline done.                                                                                               • Responsive web developer team: 1. Kimberly 2.Jared 3. Alex: start rough draft website, start looking at angular JS.                                                                                                                                                                       • Central webserver/DB developers: 1.Alex 2.Jared 3.Kimberly: rough layout for DB and start reading • Android App Developers: 1.Colton 2.Jake: Get software ready, start reading.     "
$serverText = " Get server up and running, Basic App running on an android device. "
$websiteText = "Website pull and push from DB and app pull and push from DB. "

# Assign in the same order the strings were first introduced in the
# workbook (so shared-string indices line up): "Wrap up..." first, then
# the three sprint-1 milestone cells.
$ws.Range("C13").Value = $wrapUp
$ws.Range("C2").Value = $pmText
$ws.Range("C3").Value = $serverText
$ws.Range("C4").Value = $websiteText

# Wrap text across header row + all data rows/columns A1:C13
$ws.Range("A1:C13").WrapText = $true

# Increase row height on the data rows (2-13) to fit wrapped text
$ws.Range("A2:A13").RowHeight = 99.95

# Extend column D onward formatting (wrap text) to match style used for
# the rest of the sheet.
$ws.Range("D1:XFD1").EntireColumn.WrapText = $true

# Update the view: scroll so column C is the leftmost visible column and
# select D4.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D4").Select()
